$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("I2").Value = 0.03058954393770857
$ws.Range("J2").Value = 0.01312910284463895
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 0.001284521515735388
$ws.Range("I3").Value = 0.05728587319243589
$ws.Range("J3").Value = 0.05579868708971549
$ws.Range("B4").Value = 0.08165057067603174
$ws.Range("D4").Value = 0.04327485380116956
$ws.Range("E4").Value = 0.1130204890945138
$ws.Range("G4").Value = 0.01317715959004392
$ws.Range("K4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.0003996802557953637
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 0.3481053307642928
$ws.Range("I5").Value = 0.2280311457174627
$ws.Range("J5").Value = 0.2308533916849024
$ws.Range("K5").Value = 0.02757158006362672
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("I6").Value = 0.02169076751946607
$ws.Range("J6").Value = 0.006564551422319475
$ws.Range("B7").Value = 0.04038630377524143
$ws.Range("C8").Value = 0.4368505195843259
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 0.1394765148798852
$ws.Range("H8").Value = 0.3179190751445101
$ws.Range("I8").Value = 0.06507230255839802
$ws.Range("J8").Value = 0.135667396061269
$ws.Range("B9").Value = 0.03248463564530289
$ws.Range("K9").Value = 0.003181336161187699
$ws.Range("B10").Value = 0.01053555750658472
$ws.Range("D10").Value = 0.1286549707602336
$ws.Range("E10").Value = 0.09517514871116955
$ws.Range("G10").Value = 0.1207906295754026
$ws.Range("K10").Value = 0
$ws.Range("B12").Value = 0.03511852502194907
$ws.Range("E12").Value = 0.01586252478519498
$ws.Range("K12").Value = 0.007423117709437964
$ws.Range("D13").Value = 0
$ws.Range("I13").Value = 0.0489432703003336
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("H14").Value = 0.04110468850353238
$ws.Range("I14").Value = 0.1368186874304777
$ws.Range("J14").Value = 0.3654266958424501
$ws.Range("D15").Value = 0.03333333333333335
$ws.Range("G15").Value = 0
$ws.Range("B16").Value = 0.04126426690079016
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.05287508261731649
$ws.Range("I16").Value = 0.01501668520578421
$ws.Range("K16").Value = 0.2926829268292691
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("I17").Value = 0.00389321468298109
$ws.Range("J17").Value = 0.002188183807439825
$ws.Range("B18").Value = 0.003511852502194908
$ws.Range("E18").Value = 0.001982815598149372
$ws.Range("K18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("H19").Value = 0.01220295439948619
$ws.Range("I19").Value = 0.09343715239154579
$ws.Range("J19").Value = 0.003282275711159737
$ws.Range("D20").Value = 0
$ws.Range("I20").Value = 0.02280311457174639
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("I21").Value = 0.04115684093437147
$ws.Range("J21").Value = 0.008752735229759299
$ws.Range("D22").Value = 0.001169590643274854
$ws.Range("G22").Value = 0.08345534407027819
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("I23").Value = 0.1273637374860951
$ws.Range("J23").Value = 0.0667396061269146
$ws.Range("B24").Value = 0.001755926251097454
$ws.Range("E24").Value = 0.007270323859881032
$ws.Range("K24").Value = 0
$ws.Range("B32").Value = 0.06935908691834947
$ws.Range("D32").Value = 0.0614035087719297
$ws.Range("E32").Value = 0.09781890284203537
$ws.Range("G32").Value = 0.04612005856515373
$ws.Range("K32").Value = 0.2969247083775194
$ws.Range("B33").Value = 0.05882352941176469
$ws.Range("D33").Value = 0.07602339181286531
$ws.Range("E33").Value = 0.2709847984137472
$ws.Range("G33").Value = 0
$ws.Range("K33").Value = 0.008483563096500531
$ws.Range("D34").Value = 0
$ws.Range("I34").Value = 0.003337041156840935
